$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / corrected "Samples" query (drops std.id and dgn."participant.id" columns)
$newSampleQuery = @'
SELECT DISTINCT
   smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID", std.dbgap_accession AS "Study ID", smp.anatomic_site AS "Sample Anatomic Site",
    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification",
 dgn.diagnosis as "Sample Diagnosis"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON smp."id" = dgn."sample.id"
WHERE 
   std.dbgap_accession = 'phs000468' 
AND smp.sample_id IS NOT NULL
ORDER BY 
   smp.sample_id ASC
;
'@

$ws.Range("B4").Value = $newSampleQuery
$ws.Rows.Item(4).RowHeight = 378

# Update the view: scroll to row 4 and select B4 (matches post-edit saved view state)
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Save()
